$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17=640.7816, J17=537.2958, L17=1611.8874, N17=-1947.8874
$ws.Range("H17").Value = 640.7816
$ws.Range("J17").Value = 537.2958
$ws.Range("L17").Value = 1611.8874
$ws.Range("N17").Value = -1947.8874
# Row 74: H74=4004360, I74=9094364, J74=5071.4287, K74=9094364, L74=5071.4287, M74=-9093428, N74=-6943.4287
$ws.Range("H74").Value = 4004360
$ws.Range("I74").Value = 9094364
$ws.Range("J74").Value = 5071.4287
$ws.Range("K74").Value = 9094364
$ws.Range("L74").Value = 5071.4287
$ws.Range("M74").Value = -9093428
$ws.Range("N74").Value = -6943.4287
# Row 77: H77=4004360, I77=9094364, J77=5071.4287, K77=45471820, L77=25357.1435, M77=-45467140, N77=-34717.14350000001
$ws.Range("H77").Value = 4004360
$ws.Range("I77").Value = 9094364
$ws.Range("J77").Value = 5071.4287
$ws.Range("K77").Value = 45471820
$ws.Range("L77").Value = 25357.1435
$ws.Range("M77").Value = -45467140
$ws.Range("N77").Value = -34717.14350000001
# Row 87: H87=23674, J87=23674, L87=23674, N87=-26170
$ws.Range("H87").Value = 23674
$ws.Range("J87").Value = 23674
$ws.Range("L87").Value = 23674
$ws.Range("N87").Value = -26170
# Row 90: H90=23674, J90=23674, L90=71022, N90=-83502
$ws.Range("H90").Value = 23674
$ws.Range("J90").Value = 23674
$ws.Range("L90").Value = 71022
$ws.Range("N90").Value = -83502

$ws = $wb.Worksheets.Item("ARM")
# Row 2: H2=715.375, I2=751.6667, K2=751.6667, M2=-638.6667
$ws.Range("H2").Value = 715.375
$ws.Range("I2").Value = 751.6667
$ws.Range("K2").Value = 751.6667
$ws.Range("M2").Value = -638.6667
# Row 32: H32=3544.527, I32=3220.7273, J32=6215.875, K32=3220.7273, L32=6215.875, M32=-2933.7273, N32=-6789.875
$ws.Range("H32").Value = 3544.527
$ws.Range("I32").Value = 3220.7273
$ws.Range("J32").Value = 6215.875
$ws.Range("K32").Value = 3220.7273
$ws.Range("L32").Value = 6215.875
$ws.Range("M32").Value = -2933.7273
$ws.Range("N32").Value = -6789.875
# Row 63: H63=15393113
$ws.Range("H63").Value = 15393113
# Row 66: H66=15393113
$ws.Range("H66").Value = 15393113
# Row 74: H74=2591.8076, I74=2403.6086, K74=2403.6086, M74=-1529.6086
$ws.Range("H74").Value = 2591.8076
$ws.Range("I74").Value = 2403.6086
$ws.Range("K74").Value = 2403.6086
$ws.Range("M74").Value = -1529.6086
# Row 77: H77=2591.8076, I77=2403.6086, K77=12018.043, M77=-7650.043
$ws.Range("H77").Value = 2591.8076
$ws.Range("I77").Value = 2403.6086
$ws.Range("K77").Value = 12018.043
$ws.Range("M77").Value = -7650.043
# Row 110: H110=2098, I110=2089.3125, J110=2117.8572, K110=2089.3125, L110=2117.8572, M110=-44.3125, N110=-6207.8572
$ws.Range("H110").Value = 2098
$ws.Range("I110").Value = 2089.3125
$ws.Range("J110").Value = 2117.8572
$ws.Range("K110").Value = 2089.3125
$ws.Range("L110").Value = 2117.8572
$ws.Range("M110").Value = -44.3125
$ws.Range("N110").Value = -6207.8572
# Row 116: H116=715.375, I116=751.6667, K116=751.6667, M116=1542.3333
$ws.Range("H116").Value = 715.375
$ws.Range("I116").Value = 751.6667
$ws.Range("K116").Value = 751.6667
$ws.Range("M116").Value = 1542.3333
# Row 122: H122=1707.5652, I122=1229.1945, J122=3429.7, K122=3687.5835, L122=10289.1, M122=-1237.5835, N122=-15189.1
$ws.Range("H122").Value = 1707.5652
$ws.Range("I122").Value = 1229.1945
$ws.Range("J122").Value = 3429.7
$ws.Range("K122").Value = 3687.5835
$ws.Range("L122").Value = 10289.1
$ws.Range("M122").Value = -1237.5835
$ws.Range("N122").Value = -15189.1

$ws = $wb.Worksheets.Item("BSM")
# Row 3: H3=715.375, I3=751.6667, K3=751.6667, M3=-637.6667
$ws.Range("H3").Value = 715.375
$ws.Range("I3").Value = 751.6667
$ws.Range("K3").Value = 751.6667
$ws.Range("M3").Value = -637.6667
# Row 82: H82=22788.666, J82=29714.545, L82=29714.545, N82=-30480.545
$ws.Range("H82").Value = 22788.666
$ws.Range("J82").Value = 29714.545
$ws.Range("L82").Value = 29714.545
$ws.Range("N82").Value = -30480.545
# Row 85: H85=22788.666, J85=29714.545, L85=29714.545, N85=-32366.545
$ws.Range("H85").Value = 22788.666
$ws.Range("J85").Value = 29714.545
$ws.Range("L85").Value = 29714.545
$ws.Range("N85").Value = -32366.545
# Row 107: H107=2087, I107=2337, J107=1837, K107=2337, L107=1837, M107=-417, N107=-5677
$ws.Range("H107").Value = 2087
$ws.Range("I107").Value = 2337
$ws.Range("J107").Value = 1837
$ws.Range("K107").Value = 2337
$ws.Range("L107").Value = 1837
$ws.Range("M107").Value = -417
$ws.Range("N107").Value = -5677

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31=9618198, I31=1469.9584, K31=1469.9584, M31=-1174.9584
$ws.Range("H31").Value = 9618198
$ws.Range("I31").Value = 1469.9584
$ws.Range("K31").Value = 1469.9584
$ws.Range("M31").Value = -1174.9584
# Row 34: H34=9618198, I34=1469.9584, K34=1469.9584, M34=-1267.9584
$ws.Range("H34").Value = 9618198
$ws.Range("I34").Value = 1469.9584
$ws.Range("K34").Value = 1469.9584
$ws.Range("M34").Value = -1267.9584

$ws = $wb.Worksheets.Item("CUL")
# Row 34: H34=6088.769, I34=15839.143, J34=3955.875, K34=47517.429, L34=11867.625, M34=-47433.429, N34=-12035.625
$ws.Range("H34").Value = 6088.769
$ws.Range("I34").Value = 15839.143
$ws.Range("J34").Value = 3955.875
$ws.Range("K34").Value = 47517.429
$ws.Range("L34").Value = 11867.625
$ws.Range("M34").Value = -47433.429
$ws.Range("N34").Value = -12035.625
# Row 39: H39=10791.048, J39=10930.6, L39=32791.8, N39=-33379.8
$ws.Range("H39").Value = 10791.048
$ws.Range("J39").Value = 10930.6
$ws.Range("L39").Value = 32791.8
$ws.Range("N39").Value = -33379.8
# Row 55: H55=4452.5, J55=4943, L55=14829, N55=-15183
$ws.Range("H55").Value = 4452.5
$ws.Range("J55").Value = 4943
$ws.Range("L55").Value = 14829
$ws.Range("N55").Value = -15183
# Row 87: H87=2500, I87=2500, K87=7500, M87=-6252
$ws.Range("H87").Value = 2500
$ws.Range("I87").Value = 2500
$ws.Range("K87").Value = 7500
$ws.Range("M87").Value = -6252
# Row 90: H90=2500, I90=2500, K90=22500, M90=-16260
$ws.Range("H90").Value = 2500
$ws.Range("I90").Value = 2500
$ws.Range("K90").Value = 22500
$ws.Range("M90").Value = -16260
# Row 113: H113=584.53656, I113=522.4583, K113=1567.3749, M113=602.6251
$ws.Range("H113").Value = 584.53656
$ws.Range("I113").Value = 522.4583
$ws.Range("K113").Value = 1567.3749
$ws.Range("M113").Value = 602.6251
# Row 131: H131=6098406.5, I131=71429130, J131=872.72, K131=214287390, L131=2618.16, M131=-214282350, N131=-12698.16
$ws.Range("H131").Value = 6098406.5
$ws.Range("I131").Value = 71429130
$ws.Range("J131").Value = 872.72
$ws.Range("K131").Value = 214287390
$ws.Range("L131").Value = 2618.16
$ws.Range("M131").Value = -214282350
$ws.Range("N131").Value = -12698.16
# Row 132: H132=1775.4706, I132=722.875, J132=2711.111, K132=6505.875, L132=24399.999, M132=-3975.875, N132=-29459.999
$ws.Range("H132").Value = 1775.4706
$ws.Range("I132").Value = 722.875
$ws.Range("J132").Value = 2711.111
$ws.Range("K132").Value = 6505.875
$ws.Range("L132").Value = 24399.999
$ws.Range("M132").Value = -3975.875
$ws.Range("N132").Value = -29459.999

$ws = $wb.Worksheets.Item("GSM")
# Row 107: H107=505.5909, I107=261.25, J107=798.8, K107=261.25, L107=798.8, M107=1658.75, N107=-4638.8
$ws.Range("H107").Value = 505.5909
$ws.Range("I107").Value = 261.25
$ws.Range("J107").Value = 798.8
$ws.Range("K107").Value = 261.25
$ws.Range("L107").Value = 798.8
$ws.Range("M107").Value = 1658.75
$ws.Range("N107").Value = -4638.8

$ws = $wb.Worksheets.Item("LTW")
# Row 122: H122=5416.9473, I122=2891.7, J122=8222.777, K122=8675.099999999999, L122=24668.331, M122=-6225.099999999999, N122=-29568.331
$ws.Range("H122").Value = 5416.9473
$ws.Range("I122").Value = 2891.7
$ws.Range("J122").Value = 8222.777
$ws.Range("K122").Value = 8675.099999999999
$ws.Range("L122").Value = 24668.331
$ws.Range("M122").Value = -6225.099999999999
$ws.Range("N122").Value = -29568.331

$ws = $wb.Worksheets.Item("WVR")
# Row 64: H64=25900, J64=25900, L64=25900, N64=-26396
$ws.Range("H64").Value = 25900
$ws.Range("J64").Value = 25900
$ws.Range("L64").Value = 25900
$ws.Range("N64").Value = -26396
# Row 67: H67=25900, J67=25900, L67=25900, N67=-27616
$ws.Range("H67").Value = 25900
$ws.Range("J67").Value = 25900
$ws.Range("L67").Value = 25900
$ws.Range("N67").Value = -27616
# Row 136: H136=1845.7301, I136=558.7273, J136=4826.1577, K136=1676.1819, L136=14478.4731, M136=873.8181, N136=-19578.4731
$ws.Range("H136").Value = 1845.7301
$ws.Range("I136").Value = 558.7273
$ws.Range("J136").Value = 4826.1577
$ws.Range("K136").Value = 1676.1819
$ws.Range("L136").Value = 14478.4731
$ws.Range("M136").Value = 873.8181
$ws.Range("N136").Value = -19578.4731
